$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "data" to "VytrackUsers"
$ws.Name = "VytrackUsers"

# Delete column E (Result/PASS/FAIL column) entirely
$ws.Columns.Item(5).Delete()

# Move selection to F1 (matches the post-edit sheetView selection)
$ws.Range("F1").Select()
